# Hortaliza, Macroferia Regional de Talca - Poroto verde
# A new weekly price record (row 61) is inserted; all subsequent records
# shift down by one row (old row 61 -> new row 62, ..., old row 87 -> new row 88).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, pushing existing rows 61-87 down to 62-88.
$ws.Rows.Item(61).Insert()

# Fill in the newly inserted row 61 with the new weekly record.
$ws.Cells.Item(61, 1).Value  = 5
$ws.Cells.Item(61, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(61, 3).Value  = "Maule"
$ws.Cells.Item(61, 4).Value  = 44460
$ws.Cells.Item(61, 5).Value  = 7
$ws.Cells.Item(61, 6).Value  = 100112031
$ws.Cells.Item(61, 7).Value  = "Poroto verde"
$ws.Cells.Item(61, 8).Value  = "Sin especificar"
$ws.Cells.Item(61, 9).Value  = "Primera"
$ws.Cells.Item(61, 10).Value = 200
$ws.Cells.Item(61, 11).Value = 32000
$ws.Cells.Item(61, 12).Value = 32000
$ws.Cells.Item(61, 13).Value = 32000
$ws.Cells.Item(61, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(61, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(61, 16).Value = 1280
$ws.Cells.Item(61, 17).Value = 25
$ws.Cells.Item(61, 18).Value = "Hortaliza"
